$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "Nagpur" / "Maharashtra" right before the existing
# row 143 (shifts all subsequent city rows down by one, A1:B214 -> A1:B215).
$ws.Rows("143:143").Insert()
$ws.Range("A143").Value = "Nagpur"
$ws.Range("B143").Value = "Maharashtra"

# Move the active selection from A141 to F141, matching the saved view state.
$ws.Range("F141").Select()
